$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 13.091
$ws.Range("C3").Value = -12.994
$ws.Range("E5").Value = 12.929
$ws.Range("C14").Value = -12.431
$ws.Range("C16").Value = -11.997
$ws.Range("E16").Value = 13.108
$ws.Range("C21").Value = -12.953
$ws.Range("C23").Value = -12.813
$ws.Range("C25").Value = -13.175
